# Apply the CCS Percentages by Entity edits:
#  - CPbE-FoCSbS (sheet "CPbE-FoCSbS"): row 2 becomes hardcoded 1s, row 3 B3 becomes
#    hardcoded 0 (C3:AM3 keep their formulas but recalc to 0).
#  - CPbE-FoESCbES (sheet "CPbE-FoESCbES"): rows 2-11 and 14 become hardcoded 0s,
#    row 12 becomes hardcoded 1s, row 13 formulas are redirected to reference the
#    next column over (B13=C13, C13=D13, ... AM13=AN13).

$wb = $excel.ActiveWorkbook

# ---- CPbE-FoCSbS ----
$wsA = $wb.Worksheets.Item("CPbE-FoCSbS")
$wsA.Range("B2:AM2").Value = 1
$wsA.Range("B3").Value = 0

# ---- CPbE-FoESCbES ----
$wsB = $wb.Worksheets.Item("CPbE-FoESCbES")
$wsB.Range("B2:AM11").Value = 0
$wsB.Range("B12:AM12").Value = 1
$wsB.Range("B14:AM14").Value = 0

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN")
for ($i = 0; $i -lt $cols.Length - 1; $i++) {
    $cell = $cols[$i] + "13"
    $nextcell = $cols[$i + 1] + "13"
    $wsB.Range($cell).Formula = "=" + $nextcell
}

# ---- Selections recorded in the saved file ----
$wsA.Range("B4").Select()
$wsB.Range("B14:AM14").Select()

# Restore the originally-active sheet/tab so we don't disturb unrelated state.
$wb.Worksheets.Item("About").Activate()
